$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A49").Value = "2025-04-29 05:44:03"
$ws.Range("B49").Value = 149
